$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.581137
$ws.Range("H2").Value = 1.743411
$ws.Range("I2").Value = 0.08244918404718141
$ws.Range("J2").Value = 0.08244918404718142
$ws.Range("M2").Value = 11.61289466666667
$ws.Range("N2").Value = 34.838684
$ws.Range("O2").Value = 0.09693042549509606
$ws.Range("P2").Value = 0.09693042549509606
$ws.Range("Q2").Value = 6.748682767902667
$ws.Range("R2").Value = 60.738144911124
$ws.Range("S2").Value = 0.00799183449141678
$ws.Range("T2").Value = 0.007991834491416782
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.581137
$ws.Range("H3").Value = 1.743411
$ws.Range("I3").Value = 0.08244918404718141
$ws.Range("J3").Value = 0.08244918404718142
$ws.Range("O3").Value = 0.2981108740043866
$ws.Range("P3").Value = 0.2981108740043866
$ws.Range("Q3").Value = 20.75566787251534
$ws.Range("R3").Value = 186.801010852638
$ws.Range("S3").Value = 0.02457899831725378
$ws.Range("T3").Value = 0.02457899831725378
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.581137
$ws.Range("H4").Value = 1.743411
$ws.Range("I4").Value = 0.08244918404718141
$ws.Range("J4").Value = 0.08244918404718142
$ws.Range("M4").Value = 27.39934733333333
$ws.Range("N4").Value = 82.198042
$ws.Range("O4").Value = 0.2286966748205465
$ws.Range("P4").Value = 0.2286966748205465
$ws.Range("Q4").Value = 15.92277451125134
$ws.Range("R4").Value = 143.304970601262
$ws.Range("S4").Value = 0.01885585423325764
$ws.Range("T4").Value = 0.01885585423325764
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.581137
$ws.Range("H5").Value = 1.743411
$ws.Range("I5").Value = 0.08244918404718141
$ws.Range("J5").Value = 0.08244918404718142
$ws.Range("M5").Value = 45.078635
$ws.Range("N5").Value = 135.235905
$ws.Range("O5").Value = 0.3762620256799708
$ws.Range("P5").Value = 0.3762620256799709
$ws.Range("Q5").Value = 26.196862707995
$ws.Range("R5").Value = 235.771764371955
$ws.Range("S5").Value = 0.03102249700525321
$ws.Range("T5").Value = 0.03102249700525322
$ws.Range("H6").Value = 5.486604999999999
$ws.Range("I6").Value = 0.2594718660368586
$ws.Range("J6").Value = 0.2594718660368586
$ws.Range("M6").Value = 11.61289466666667
$ws.Range("N6").Value = 34.838684
$ws.Range("O6").Value = 0.09693042549509606
$ws.Range("P6").Value = 0.09693042549509606
$ws.Range("Q6").Value = 21.23845531420222
$ws.Range("R6").Value = 191.14609782782
$ws.Range("S6").Value = 0.02515071837895927
$ws.Range("T6").Value = 0.02515071837895927
$ws.Range("H7").Value = 5.486604999999999
$ws.Range("I7").Value = 0.2594718660368586
$ws.Range("J7").Value = 0.2594718660368586
$ws.Range("O7").Value = 0.2981108740043866
$ws.Range("P7").Value = 0.2981108740043866
$ws.Range("Q7").Value = 65.31916520412111
$ws.Range("R7").Value = 587.8724868370899
$ws.Range("S7").Value = 0.07735138476379702
$ws.Range("T7").Value = 0.07735138476379703
$ws.Range("H8").Value = 5.486604999999999
$ws.Range("I8").Value = 0.2594718660368586
$ws.Range("J8").Value = 0.2594718660368586
$ws.Range("M8").Value = 27.39934733333333
$ws.Range("N8").Value = 82.198042
$ws.Range("O8").Value = 0.2286966748205465
$ws.Range("P8").Value = 0.2286966748205465
$ws.Range("Q8").Value = 50.10979869193444
$ws.Range("R8").Value = 450.9881882274099
$ws.Range("S8").Value = 0.05934035297211186
$ws.Range("T8").Value = 0.05934035297211187
$ws.Range("H9").Value = 5.486604999999999
$ws.Range("I9").Value = 0.2594718660368586
$ws.Range("J9").Value = 0.2594718660368586
$ws.Range("M9").Value = 45.078635
$ws.Range("N9").Value = 135.235905
$ws.Range("O9").Value = 0.3762620256799708
$ws.Range("P9").Value = 0.3762620256799709
$ws.Range("Q9").Value = 82.44288806139166
$ws.Range("R9").Value = 741.9859925525249
$ws.Range("S9").Value = 0.09762940992199043
$ws.Range("T9").Value = 0.09762940992199046
$ws.Range("G10").Value = 4.568238333333333
$ws.Range("H10").Value = 13.704715
$ws.Range("I10").Value = 0.6481217391361921
$ws.Range("J10").Value = 0.6481217391361921
$ws.Range("M10").Value = 11.61289466666667
$ws.Range("N10").Value = 34.838684
$ws.Range("O10").Value = 0.09693042549509606
$ws.Range("P10").Value = 0.09693042549509606
$ws.Range("Q10").Value = 53.05047057722889
$ws.Range("R10").Value = 477.45423519506
$ws.Range("S10").Value = 0.06282271594709275
$ws.Range("T10").Value = 0.06282271594709275
$ws.Range("G11").Value = 4.568238333333333
$ws.Range("H11").Value = 13.704715
$ws.Range("I11").Value = 0.6481217391361921
$ws.Range("J11").Value = 0.6481217391361921
$ws.Range("O11").Value = 0.2981108740043866
$ws.Range("P11").Value = 0.2981108740043866
$ws.Range("Q11").Value = 163.1574613372745
$ws.Range("R11").Value = 1468.41715203547
$ws.Range("S11").Value = 0.1932121381151333
$ws.Range("T11").Value = 0.1932121381151333
$ws.Range("G12").Value = 4.568238333333333
$ws.Range("H12").Value = 13.704715
$ws.Range("I12").Value = 0.6481217391361921
$ws.Range("J12").Value = 0.6481217391361921
$ws.Range("M12").Value = 27.39934733333333
$ws.Range("N12").Value = 82.198042
$ws.Range("O12").Value = 0.2286966748205465
$ws.Range("P12").Value = 0.2286966748205465
$ws.Range("Q12").Value = 125.1667487964478
$ws.Range("R12").Value = 1126.50073916803
$ws.Range("S12").Value = 0.1482232866193568
$ws.Range("T12").Value = 0.1482232866193568
$ws.Range("G13").Value = 4.568238333333333
$ws.Range("H13").Value = 13.704715
$ws.Range("I13").Value = 0.6481217391361921
$ws.Range("J13").Value = 0.6481217391361921
$ws.Range("M13").Value = 45.078635
$ws.Range("N13").Value = 135.235905
$ws.Range("O13").Value = 0.3762620256799708
$ws.Range("P13").Value = 0.3762620256799709
$ws.Range("Q13").Value = 205.9299484213417
$ws.Range("R13").Value = 1853.369535792075
$ws.Range("S13").Value = 0.2438635984546093
$ws.Range("T13").Value = 0.2438635984546093
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.07018266666666667
$ws.Range("H14").Value = 0.210548
$ws.Range("I14").Value = 0.009957210779767909
$ws.Range("J14").Value = 0.009957210779767911
$ws.Range("M14").Value = 11.61289466666667
$ws.Range("N14").Value = 34.838684
$ws.Range("O14").Value = 0.09693042549509606
$ws.Range("P14").Value = 0.09693042549509606
$ws.Range("Q14").Value = 0.8150239154257779
$ws.Range("R14").Value = 7.335215238832
$ws.Range("S14").Value = 0.0009651566776272607
$ws.Range("T14").Value = 0.0009651566776272608
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.07018266666666667
$ws.Range("H15").Value = 0.210548
$ws.Range("I15").Value = 0.009957210779767909
$ws.Range("J15").Value = 0.009957210779767911
$ws.Range("O15").Value = 0.2981108740043866
$ws.Range("P15").Value = 0.2981108740043866
$ws.Range("Q15").Value = 2.506617406464889
$ws.Range("R15").Value = 22.559556658184
$ws.Range("S15").Value = 0.002968352808202511
$ws.Range("T15").Value = 0.002968352808202511
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.07018266666666667
$ws.Range("H16").Value = 0.210548
$ws.Range("I16").Value = 0.009957210779767909
$ws.Range("J16").Value = 0.009957210779767911
$ws.Range("M16").Value = 27.39934733333333
$ws.Range("N16").Value = 82.198042
$ws.Range("O16").Value = 0.2286966748205465
$ws.Range("P16").Value = 0.2286966748205465
$ws.Range("Q16").Value = 1.922959260779556
$ws.Range("R16").Value = 17.306633347016
$ws.Range("S16").Value = 0.002277180995820222
$ws.Range("T16").Value = 0.002277180995820222
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.07018266666666667
$ws.Range("H17").Value = 0.210548
$ws.Range("I17").Value = 0.009957210779767909
$ws.Range("J17").Value = 0.009957210779767911
$ws.Range("M17").Value = 45.078635
$ws.Range("N17").Value = 135.235905
$ws.Range("O17").Value = 0.3762620256799708
$ws.Range("P17").Value = 0.3762620256799709
$ws.Range("Q17").Value = 3.163738813993334
$ws.Range("R17").Value = 28.47364932594
$ws.Range("S17").Value = 0.003746520298117915
$ws.Range("T17").Value = 0.003746520298117916
